# Add Sections to the presentation (Slide Sorter "Sections" feature).
#
# Target layout (slide ids from the underlying sldIdLst):
#   Default Section     -> slide 256            (Slides.Item(1))
#   Machine Learning     -> slides 257-264        (Slides.Item(2)  .. Item(9))
#   Linear Regression    -> slides 265,266,268    (Slides.Item(10) .. Item(12))

$p = $ppt.ActivePresentation
$sp = $p.SectionProperties

# AddBeforeSlide(SlideIndex, Name) inserts a new section that starts at
# SlideIndex and runs through the slide before the next section's start
# (or through the end of the deck for the last section added).
$sp.AddBeforeSlide(1, "Default Section") | Out-Null
$sp.AddBeforeSlide(2, "Machine Learning") | Out-Null
$sp.AddBeforeSlide(10, "Linear Regression") | Out-Null
